$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Front End")
$ws2 = $wb.Worksheets.Item("Back End")

# ---------------------------------------------------------------------------
# Front End sheet: mark several in-flight tasks complete, add a new task row
# ---------------------------------------------------------------------------

# Code 200 message -> COMPLETE (with completion date)
$ws1.Range("B4").Value = "COMPLETE"
$ws1.Range("C3").Copy()
$ws1.Range("C4").PasteSpecial(-4122)
$ws1.Range("C4").Value = "6/2/2023"

# Info display from JS object -> COMPLETE (with completion date)
$ws1.Range("B5").Value = "COMPLETE"
$ws1.Range("C3").Copy()
$ws1.Range("C5").PasteSpecial(-4122)
$ws1.Range("C5").Value = "6/2/2023"

# New task: Data table -> COMPLETE (with completion date)
$ws1.Range("A7").Value = "Data table"
$ws1.Range("B7").Value = "COMPLETE"
$ws1.Range("C3").Copy()
$ws1.Range("C7").PasteSpecial(-4122)
$ws1.Range("C7").Value = "6/2/2023"

# Site name -> COMPLETE (with completion date)
$ws1.Range("B15").Value = "COMPLETE"
$ws1.Range("C3").Copy()
$ws1.Range("C15").PasteSpecial(-4122)
$ws1.Range("C15").Value = "6/2/2023"

# ---------------------------------------------------------------------------
# Back End sheet: mark several in-flight tasks complete, add two new ones
# ---------------------------------------------------------------------------

# Searching via service tag -> COMPLETE (with completion date)
$ws2.Range("B10").Value = "COMPLETE"
$ws2.Range("C3").Copy()
$ws2.Range("C10").PasteSpecial(-4122)
$ws2.Range("C10").Value = "6/2/2023"

# Write to CSV -> COMPLETE (with completion date)
$ws2.Range("B17").Value = "COMPLETE"
$ws2.Range("C3").Copy()
$ws2.Range("C17").PasteSpecial(-4122)
$ws2.Range("C17").Value = "6/2/2023"

# Mass-write to CSV -> COMPLETE (with completion date)
$ws2.Range("B18").Value = "COMPLETE"
$ws2.Range("C3").Copy()
$ws2.Range("C18").PasteSpecial(-4122)
$ws2.Range("C18").Value = "6/2/2023"

# New conditional formatting (Highlight Cell Rules style) for the two newly
# populated status cells, matching the IN PROGRESS / COMPLETE / INCOMPLETE
# rules already present on the rest of the Status column.
$rngB17 = $ws2.Range("B17")
$fc = $rngB17.FormatConditions.Add(1, 3, '"IN PROGRESS"')
$fc.Font.Color = 22428
$fc.Interior.Color = 10284031
$fc = $rngB17.FormatConditions.Add(1, 3, '"COMPLETE"')
$fc.Font.Color = 24832
$fc.Interior.Color = 13561798
$fc = $rngB17.FormatConditions.Add(1, 3, '"INCOMPLETE"')
$fc.Font.Color = 393372
$fc.Interior.Color = 13551615

$rngB18 = $ws2.Range("B18")
$fc = $rngB18.FormatConditions.Add(1, 3, '"IN PROGRESS"')
$fc.Font.Color = 22428
$fc.Interior.Color = 10284031
$fc = $rngB18.FormatConditions.Add(1, 3, '"COMPLETE"')
$fc.Font.Color = 24832
$fc.Interior.Color = 13561798
$fc = $rngB18.FormatConditions.Add(1, 3, '"INCOMPLETE"')
$fc.Font.Color = 393372
$fc.Interior.Color = 13551615

# ---------------------------------------------------------------------------
# Selection / active sheet: user ends up on the Back End tab
# ---------------------------------------------------------------------------

$ws1.Range("D15").Select()
$ws2.Range("C19").Select()
$ws2.Activate()
